$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the three stored SQL queries in B2:B4 to drop the trailing
# "LIMIT 100" clause (kept as a bare terminating semicolon). ---

# B4: FilesTab query - "LIMIT 100" clause (and its preceding newline) is
# removed outright; the semicolon attaches directly to "ASC".
$b4 = $ws.Range("B4").Value2
$b4 = $b4.Replace("    f1.file_name ASC`nLIMIT 100;", "    f1.file_name ASC;")
$ws.Range("B4").Value2 = $b4

# B3: SamplesTab query - "LIMIT 100;" becomes just ";" on its own line.
$b3 = $ws.Range("B3").Value2
$b3 = $b3.Replace("    smp.sample_id ASC`nLIMIT 100;", "    smp.sample_id ASC`n;")
$ws.Range("B3").Value2 = $b3

# B2: ParticipantsTab query - "LIMIT 100;" becomes just ";" on its own line.
$b2 = $ws.Range("B2").Value2
$b2 = $b2.Replace("    gender`nLIMIT 100;", "    gender`n;")
$ws.Range("B2").Value2 = $b2

# Re-assigning B2's (very long, many-line) text makes the headless engine's
# row autofit compute a height taller than Excel's 409.5pt row-height cap;
# real Excel clamps to that cap (the row's height is unaffected by this
# edit), so restore it explicitly here.
$ws.Rows(2).RowHeight = 409.5

# --- Move the active selection from C2 to C3 ---
$ws.Range("C3").Select()
